$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new record was logged between the existing 2026/01/31 row (row 727) and
# the 2026/12/29 row (old row 728): insert a fresh row at 728 and push
# everything below it down by one (old 728..769 -> new 729..770).
$ws.Rows(728).EntireRow.Insert()

# Column A holds plain text dates (e.g. "2026/01/31"), not real date
# serials. Assigning that literal string directly gets auto-coerced to a
# date value by the COM layer, so stage it through a text number format,
# then clear the format again so the cell matches its neighbours (no
# explicit style index).
$ws.Range("A728").NumberFormat = "@"
$ws.Range("A728").Value = "2026/01/31"
$ws.Range("B728").Value = "土"
$ws.Range("C728").Value = 8
$ws.Range("D728").Value = 25
$ws.Range("A728").ClearFormats()
